$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for columns G:T, rows 2-10
$updates = @{
    "G2" = 12.31250333333333
    "H2" = 36.93751
    "I2" = 0.6498350963072504
    "J2" = 0.6498350963072506
    "K2" = 2
    "L2" = 0.6666666666666666
    "M2" = 0.07259900000000001
    "N2" = 0.217797
    "O2" = 0.0162094769588191
    "P2" = 0.0162094769588191
    "Q2" = 0.8938754294966669
    "R2" = 8.044878865470002
    "S2" = 0.01053348702062437
    "T2" = 0.01053348702062437
    "G3" = 12.31250333333333
    "H3" = 36.93751
    "I3" = 0.6498350963072504
    "J3" = 0.6498350963072506
    "O3" = 0.9349228167457665
    "P3" = 0.9349228167457664
    "Q3" = 51.55653920777334
    "R3" = 464.00885286996
    "S3" = 0.607545658659831
    "T3" = 0.607545658659831
    "G4" = 12.31250333333333
    "H4" = 36.93751
    "I4" = 0.6498350963072504
    "J4" = 0.6498350963072506
    "M4" = 0.2188686666666667
    "N4" = 0.656606
    "O4" = 0.04886770629541442
    "P4" = 0.04886770629541441
    "Q4" = 2.694821187895556
    "R4" = 24.25339069106
    "S4" = 0.03175595062679506
    "T4" = 0.03175595062679506
    "I5" = 0.3333514949915254
    "J5" = 0.3333514949915254
    "K5" = 2
    "L5" = 0.6666666666666666
    "M5" = 0.07259900000000001
    "N5" = 0.217797
    "O5" = 0.0162094769588191
    "P5" = 0.0162094769588191
    "Q5" = 0.4585389623493334
    "R5" = 4.126850661144
    "S5" = 0.005403453377253032
    "T5" = 0.005403453377253033
    "I6" = 0.3333514949915254
    "J6" = 0.3333514949915254
    "O6" = 0.9349228167457665
    "P6" = 0.9349228167457664
    "S6" = 0.3116579186638892
    "T6" = 0.3116579186638892
    "I7" = 0.3333514949915254
    "J7" = 0.3333514949915254
    "M7" = 0.2188686666666667
    "N7" = 0.656606
    "O7" = 0.04886770629541442
    "P7" = 0.04886770629541441
    "Q7" = 1.382385588012444
    "R7" = 12.441470292112
    "S7" = 0.01629012295038317
    "T7" = 0.01629012295038317
    "G8" = 0.3185656666666667
    "H8" = 0.955697
    "I8" = 0.01681340870122405
    "J8" = 0.01681340870122405
    "K8" = 2
    "L8" = 0.6666666666666666
    "M8" = 0.07259900000000001
    "N8" = 0.217797
    "O8" = 0.0162094769588191
    "P8" = 0.0162094769588191
    "Q8" = 0.02312754883433334
    "R8" = 0.208147939509
    "S8" = 0.0002725365609416998
    "T8" = 0.0002725365609416998
    "G9" = 0.3185656666666667
    "H9" = 0.955697
    "I9" = 0.01681340870122405
    "J9" = 0.01681340870122405
    "O9" = 0.9349228167457665
    "P9" = 0.9349228167457664
    "Q9" = 1.333940210134667
    "R9" = 12.005461891212
    "S9" = 0.01571923942204617
    "T9" = 0.01571923942204616
    "G10" = 0.3185656666666667
    "H10" = 0.955697
    "I10" = 0.01681340870122405
    "J10" = 0.01681340870122405
    "M10" = 0.2188686666666667
    "N10" = 0.656606
    "O10" = 0.04886770629541442
    "P10" = 0.04886770629541441
    "Q10" = 0.06972404270911112
    "R10" = 0.6275163843820001
    "S10" = 0.0008216327182361821
    "T10" = 0.0008216327182361819
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
